# Generate Report for Handoff
# Update the status/priority/handoff-datetime for the 779dd645-... file
# from "In Translation" to "Ready for handoff" across the Overview, zh-cn
# and de-de sheets, plus the related Priority/Latest Handoff Datetime values.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-12 04:14:32"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-12 04:14:27"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-12 04:14:32"

# Column widths widen slightly to fit the new longer "Ready for handoff" text
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 16.3
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 16.3
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 16.3
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 16.3
